$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (D1:E1) --------------------------------------------
$ws.Range("D1").Value = "Baseweight"
$ws.Range("E1").Value = "Lengthweight"

# --- Baseweight / Lengthweight data (D2:E33) ------------------------------
$ws.Range("D2").Value = 5.9
$ws.Range("E2").Value = 0.16
$ws.Range("D3").Value = 23.9
$ws.Range("E3").Value = 0.24545500000000001
$ws.Range("D4").Value = 23.9
$ws.Range("E4").Value = 0.22272700000000001
$ws.Range("D5").Value = 23.9
$ws.Range("E5").Value = 0.21363599999999999
$ws.Range("D6").Value = 26.9
$ws.Range("E6").Value = 0.24545454545454543
$ws.Range("D7").Value = 26.9
$ws.Range("E7").Value = 0.22272727272727275
$ws.Range("D8").Value = 26.9
$ws.Range("E8").Value = 0.21363636363636362
$ws.Range("D9").Value = 31.9
$ws.Range("E9").Value = 0.24545454545454543
$ws.Range("D10").Value = 31.9
$ws.Range("E10").Value = 0.22272727272727275
$ws.Range("D11").Value = 31.9
$ws.Range("E11").Value = 0.21363636363636362
$ws.Range("D12").Value = 60.4
$ws.Range("E12").Value = 0.58181818181818179
$ws.Range("D13").Value = 60.4
$ws.Range("E13").Value = 0.53181818181818186
$ws.Range("D14").Value = 60.4
$ws.Range("E14").Value = 0.50454545454545452
$ws.Range("D15").Value = 61.4
$ws.Range("E15").Value = 0.58181818181818179
$ws.Range("D16").Value = 61.4
$ws.Range("E16").Value = 0.53181818181818186
$ws.Range("D17").Value = 61.4
$ws.Range("E17").Value = 0.50454545454545452
$ws.Range("D18").Value = 81.400000000000006
$ws.Range("E18").Value = 0.58181818181818179
$ws.Range("D19").Value = 81.400000000000006
$ws.Range("E19").Value = 0.53181818181818186
$ws.Range("D20").Value = 81.400000000000006
$ws.Range("E20").Value = 0.50454545454545452
$ws.Range("D21").Value = 106.7
$ws.Range("E21").Value = 0.80909090909090908
$ws.Range("D22").Value = 106.7
$ws.Range("E22").Value = 0.74090909090909096
$ws.Range("D23").Value = 106.7
$ws.Range("E23").Value = 0.70000000000000007
$ws.Range("D24").Value = 109.7
$ws.Range("E24").Value = 0.80909090909090908
$ws.Range("D25").Value = 109.7
$ws.Range("E25").Value = 0.74090909090909096
$ws.Range("D26").Value = 109.7
$ws.Range("E26").Value = 0.70000000000000007
$ws.Range("D27").Value = 179.7
$ws.Range("E27").Value = 0.80909090909090908
$ws.Range("D28").Value = 179.7
$ws.Range("E28").Value = 0.74090909090909096
$ws.Range("D29").Value = 179.7
$ws.Range("E29").Value = 0.70000000000000007
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0

# --- Baseweight-only rows (KS12P / KS20P / KS25P) -------------------------
$ws.Range("D38").Value = 8
$ws.Range("D39").Value = 23
$ws.Range("D40").Value = 30

# --- Manually-typed cells (E3:E5) picked up Excel's "new cell" font -------
# (11pt black Calibri, instead of the sheet's default 12pt theme-coloured
# Calibri) - reproduce that formatting.
$rng = $ws.Range("E3:E5")
$rng.Font.Size = 11
$rng.Font.Color = 0

# --- Column E sizing (auto-fit to the new "Lengthweight" header) ---------
$ws.Columns.Item(5).ColumnWidth = 11.5

# --- Selection / scroll position ------------------------------------------
$null = $ws.Range("L43").Select()
